$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 05:20"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 245088
$ws.Range("C4").Value = 211
$ws.Range("E4").Value = 228610

# Row 21: Israel -> Israel
$ws.Range("F21").Value = 108

# Row 23: Australia -> Australia
$ws.Range("E23").Value = 4703
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 26

# Row 24: Noruega -> Noruega
$ws.Range("B24").Value = 5192
$ws.Range("C24").Value = 45
$ws.Range("E24").Value = 5110

# Row 53: Singapur -> Singapur
$ws.Range("E53").Value = 778
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 5

# Row 98: Honduras -> Honduras
$ws.Range("B98").Value = 222
$ws.Range("C98").Value = 3
$ws.Range("E98").Value = 204
$ws.Range("F98").Value = 10
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 15

# Row 185: Angola -> Republica de Africa Central
$ws.Range("A185").Value = "Republica de Africa Central"
$ws.Range("C185").Value = 5
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 8
$ws.Range("H185").Value = 0

# Row 186: Sudan -> Angola
$ws.Range("A186").Value = "Angola"
$ws.Range("D186").Value = 1
$ws.Range("E186").Value = 5

# Row 187: Fiyi -> Sudan
$ws.Range("A187").Value = "Sudan"
$ws.Range("B187").Value = 8
$ws.Range("D187").Value = 2
$ws.Range("E187").Value = 4
$ws.Range("H187").Value = 2

# Row 189: Liberia -> Fiyi
$ws.Range("A189").Value = "Fiyi"
$ws.Range("B189").Value = 7
$ws.Range("E189").Value = 7

# Row 190: Nepal -> Liberia
$ws.Range("A190").Value = "Liberia"
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 6

# Row 192: San Bartolome -> Nepal
$ws.Range("A192").Value = "Nepal"

# Row 193: Mauritania -> San Bartolome
$ws.Range("A193").Value = "San Bartolome"
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 5
$ws.Range("H193").Value = 0

# Row 194: Islas Turcas y Caicos -> Mauritania
$ws.Range("A194").Value = "Mauritania"
$ws.Range("B194").Value = 6
$ws.Range("D194").Value = 2
$ws.Range("E194").Value = 3
$ws.Range("H194").Value = 1

# Row 195: Nicaragua -> Islas Turcas y Caicos
$ws.Range("A195").Value = "Islas Turcas y Caicos"
$ws.Range("E195").Value = 5
$ws.Range("H195").Value = 0

# Row 197: Butan -> Nicaragua
$ws.Range("A197").Value = "Nicaragua"
$ws.Range("D197").Value = 0
$ws.Range("H197").Value = 1

# Row 198: Botsuana -> Butan
$ws.Range("A198").Value = "Butan"
$ws.Range("B198").Value = 5
$ws.Range("D198").Value = 1
$ws.Range("E198").Value = 4
$ws.Range("H198").Value = 0

# Row 199: Gambia -> Botsuana
$ws.Range("A199").Value = "Botsuana"
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 3

# Row 200: Malaui -> Gambia
$ws.Range("A200").Value = "Gambia"
$ws.Range("B200").Value = 4
$ws.Range("D200").Value = 2
$ws.Range("E200").Value = 1
$ws.Range("H200").Value = 1

# Row 201: Republica de Africa Central -> Malaui
$ws.Range("A201").Value = "Malaui"
